$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 9199.888999999999
$ws.Range("I8").Value = 4685.5713
$ws.Range("K8").Value = 14056.7139
$ws.Range("M8").Value = -13917.7139
$ws.Range("H17").Value = 32918.29
$ws.Range("J17").Value = 32918.29
$ws.Range("L17").Value = 98754.87
$ws.Range("N17").Value = -99090.87
$ws.Range("H19").Value = 830.7692
$ws.Range("I19").Value = 783.3333
$ws.Range("J19").Value = 871.4286
$ws.Range("K19").Value = 783.3333
$ws.Range("L19").Value = 871.4286
$ws.Range("M19").Value = -608.3333
$ws.Range("N19").Value = -1221.4286
$ws.Range("H52").Value = 3257.5
$ws.Range("J52").Value = 3257.5
$ws.Range("L52").Value = 9772.5
$ws.Range("N52").Value = -10092.5
$ws.Range("H121").Value = 1565.8889
$ws.Range("I121").Value = 942.5
$ws.Range("J121").Value = 1744
$ws.Range("K121").Value = 2827.5
$ws.Range("L121").Value = 5232
$ws.Range("M121").Value = -1080.5
$ws.Range("N121").Value = -8726
$ws.Range("H137").Value = 1225.2858
$ws.Range("I137").Value = 904.90625
$ws.Range("J137").Value = 2250.5
$ws.Range("K137").Value = 2714.71875
$ws.Range("L137").Value = 6751.5
$ws.Range("M137").Value = -164.71875
$ws.Range("N137").Value = -11851.5
$ws.Range("H138").Value = 2521.0254
$ws.Range("I138").Value = 1487.44
$ws.Range("J138").Value = 4303.069
$ws.Range("K138").Value = 4462.32
$ws.Range("L138").Value = 12909.207
$ws.Range("M138").Value = 677.6800000000003
$ws.Range("N138").Value = -23189.207

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 68066.664
$ws.Range("I2").Value = 1332.625
$ws.Range("J2").Value = 144334.14
$ws.Range("K2").Value = 1332.625
$ws.Range("L2").Value = 144334.14
$ws.Range("M2").Value = -1219.625
$ws.Range("N2").Value = -144560.14
$ws.Range("H61").Value = 3934.4666
$ws.Range("I61").Value = 2232.077
$ws.Range("K61").Value = 2232.077
$ws.Range("M61").Value = -2020.077
$ws.Range("H74").Value = 1170.4242
$ws.Range("I74").Value = 1329.0435
$ws.Range("J74").Value = 805.6
$ws.Range("K74").Value = 1329.0435
$ws.Range("L74").Value = 805.6
$ws.Range("M74").Value = -455.0435
$ws.Range("N74").Value = -2553.6
$ws.Range("H77").Value = 1170.4242
$ws.Range("I77").Value = 1329.0435
$ws.Range("J77").Value = 805.6
$ws.Range("K77").Value = 6645.2175
$ws.Range("L77").Value = 4028
$ws.Range("M77").Value = -2277.2175
$ws.Range("N77").Value = -12764
$ws.Range("H102").Value = 63804.812
$ws.Range("I102").Value = 1376.091
$ws.Range("J102").Value = 201148
$ws.Range("K102").Value = 1376.091
$ws.Range("L102").Value = 201148
$ws.Range("M102").Value = 245.9090000000001
$ws.Range("N102").Value = -204392
$ws.Range("H110").Value = 1322.7307
$ws.Range("I110").Value = 1328.9
$ws.Range("J110").Value = 1302.1666
$ws.Range("K110").Value = 1328.9
$ws.Range("L110").Value = 1302.1666
$ws.Range("M110").Value = 716.0999999999999
$ws.Range("N110").Value = -5392.1666
$ws.Range("H116").Value = 68066.664
$ws.Range("I116").Value = 1332.625
$ws.Range("J116").Value = 144334.14
$ws.Range("K116").Value = 1332.625
$ws.Range("L116").Value = 144334.14
$ws.Range("M116").Value = 961.375
$ws.Range("N116").Value = -148922.14
$ws.Range("H131").Value = 24500
$ws.Range("J131").Value = 24500
$ws.Range("L131").Value = 24500
$ws.Range("N131").Value = -34580
$ws.Range("H132").Value = 2109.426
$ws.Range("I132").Value = 1724.3513
$ws.Range("J132").Value = 2947.5293
$ws.Range("K132").Value = 5173.0539
$ws.Range("L132").Value = 8842.5879
$ws.Range("M132").Value = -2643.0539
$ws.Range("N132").Value = -13902.5879
$ws.Range("H136").Value = 3934.4666
$ws.Range("I136").Value = 2232.077
$ws.Range("K136").Value = 6696.231000000001
$ws.Range("M136").Value = -4146.231000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 68066.664
$ws.Range("I3").Value = 1332.625
$ws.Range("J3").Value = 144334.14
$ws.Range("K3").Value = 1332.625
$ws.Range("L3").Value = 144334.14
$ws.Range("M3").Value = -1218.625
$ws.Range("N3").Value = -144562.14
$ws.Range("H99").Value = 2267.8
$ws.Range("I99").Value = 1300
$ws.Range("K99").Value = 1300
$ws.Range("M99").Value = 198
$ws.Range("H105").Value = 2881.75
$ws.Range("I105").Value = 2815.923
$ws.Range("K105").Value = 2815.923
$ws.Range("M105").Value = -1068.923
$ws.Range("H122").Value = 29800
$ws.Range("J122").Value = 29800
$ws.Range("L122").Value = 29800
$ws.Range("N122").Value = -39600
$ws.Range("H134").Value = 2098.6545
$ws.Range("I134").Value = 1763.5
$ws.Range("K134").Value = 5290.5
$ws.Range("M134").Value = -2755.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1698.9615
$ws.Range("I31").Value = 1334.6364
$ws.Range("K31").Value = 1334.6364
$ws.Range("M31").Value = -1039.6364
$ws.Range("H34").Value = 1698.9615
$ws.Range("I34").Value = 1334.6364
$ws.Range("K34").Value = 1334.6364
$ws.Range("M34").Value = -1132.6364

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2483
$ws.Range("I5").Value = 2592.889
$ws.Range("K5").Value = 7778.667
$ws.Range("M5").Value = -7666.667
$ws.Range("H86").Value = 279.33334
$ws.Range("I86").Value = 294
$ws.Range("J86").Value = 250
$ws.Range("K86").Value = 882
$ws.Range("L86").Value = 750
$ws.Range("M86").Value = 304
$ws.Range("N86").Value = -3122
$ws.Range("H89").Value = 279.33334
$ws.Range("I89").Value = 294
$ws.Range("J89").Value = 250
$ws.Range("K89").Value = 2646
$ws.Range("L89").Value = 2250
$ws.Range("M89").Value = 3282
$ws.Range("N89").Value = -14106
$ws.Range("H122").Value = 745.6667
$ws.Range("I122").Value = 497.25
$ws.Range("J122").Value = 869.875
$ws.Range("K122").Value = 4475.25
$ws.Range("L122").Value = 7828.875
$ws.Range("M122").Value = -2025.25
$ws.Range("N122").Value = -12728.875
$ws.Range("H131").Value = 965.38
$ws.Range("I131").Value = 372.22223
$ws.Range("J131").Value = 1024.044
$ws.Range("K131").Value = 1116.66669
$ws.Range("L131").Value = 3072.132000000001
$ws.Range("M131").Value = 3923.33331
$ws.Range("N131").Value = -13152.132
$ws.Range("H135").Value = 2483
$ws.Range("I135").Value = 2592.889
$ws.Range("K135").Value = 23336.001
$ws.Range("M135").Value = -20801.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 32000
$ws.Range("J63").Value = 32000
$ws.Range("L63").Value = 32000
$ws.Range("N63").Value = -33372
$ws.Range("H66").Value = 32000
$ws.Range("J66").Value = 32000
$ws.Range("L66").Value = 96000
$ws.Range("N66").Value = -102864
$ws.Range("H68").Value = 28998.375
$ws.Range("I68").Value = 26996.75
$ws.Range("J68").Value = 31000
$ws.Range("K68").Value = 26996.75
$ws.Range("L68").Value = 31000
$ws.Range("M68").Value = -26185.75
$ws.Range("N68").Value = -32622
$ws.Range("H71").Value = 28998.375
$ws.Range("I71").Value = 26996.75
$ws.Range("J71").Value = 31000
$ws.Range("K71").Value = 80990.25
$ws.Range("L71").Value = 93000
$ws.Range("M71").Value = -76934.25
$ws.Range("N71").Value = -101112
$ws.Range("H80").Value = 3722.7778
$ws.Range("I80").Value = 3875.625
$ws.Range("K80").Value = 3875.625
$ws.Range("M80").Value = -2877.625
$ws.Range("H83").Value = 3722.7778
$ws.Range("I83").Value = 3875.625
$ws.Range("K83").Value = 19378.125
$ws.Range("M83").Value = -14386.125
$ws.Range("H97").Value = 168522.11
$ws.Range("I97").Value = 64587.375
$ws.Range("K97").Value = 64587.375
$ws.Range("M97").Value = -64091.375
$ws.Range("H102").Value = 3014.1853
$ws.Range("I102").Value = 2646.1765
$ws.Range("J102").Value = 3639.8
$ws.Range("K102").Value = 2646.1765
$ws.Range("L102").Value = 3639.8
$ws.Range("M102").Value = -1024.1765
$ws.Range("N102").Value = -6883.8
$ws.Range("H132").Value = 1421.55
$ws.Range("I132").Value = 1011.7619
$ws.Range("J132").Value = 2377.7222
$ws.Range("K132").Value = 3035.2857
$ws.Range("L132").Value = 7133.1666
$ws.Range("M132").Value = -505.2856999999999
$ws.Range("N132").Value = -12193.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 820
$ws.Range("I22").Value = 605
$ws.Range("J22").Value = 1250
$ws.Range("K22").Value = 605
$ws.Range("L22").Value = 1250
$ws.Range("M22").Value = -310
$ws.Range("N22").Value = -1840
$ws.Range("H27").Value = 820
$ws.Range("I27").Value = 605
$ws.Range("J27").Value = 1250
$ws.Range("K27").Value = 605
$ws.Range("L27").Value = 1250
$ws.Range("M27").Value = -498
$ws.Range("N27").Value = -1464
$ws.Range("H46").Value = 968.75
$ws.Range("I46").Value = 708.3333
$ws.Range("J46").Value = 1750
$ws.Range("K46").Value = 708.3333
$ws.Range("L46").Value = 1750
$ws.Range("M46").Value = -520.3333
$ws.Range("N46").Value = -2126

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 23114
$ws.Range("J64").Value = 23114
$ws.Range("L64").Value = 23114
$ws.Range("N64").Value = -23610
$ws.Range("H67").Value = 23114
$ws.Range("J67").Value = 23114
$ws.Range("L67").Value = 23114
$ws.Range("N67").Value = -24830
